$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column D ("cartQuery"), shifting old D (dbExcel) -> E and old E (WebExcel) -> F
$ws.Columns.Item(4).Insert()
$ws.Columns.Item(4).ColumnWidth = 74.8

# Header for the new column
$ws.Range("D1").Value2 = "cartQuery"

# New query text (Cypher "cart" query) for rows 2-4, same text in each row
$cartQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
  WHERE demo.breed IN ['Yorkshire Terrier']
MATCH (f:file)-[*]->(c)
WITH COLLECT(DISTINCT f.uuid) AS uuids
MATCH (f:file)
  WHERE f.uuid in uuids
OPTIONAL MATCH (f)-->(parent)
OPTIONAL MATCH (f)-[*]->(samp:sample)
OPTIONAL MATCH (f:file)-[*]->(c:case)
OPTIONAL MATCH (s:study)<-[:member_of]-(c)
OPTIONAL MATCH (c)-->(i:canine_individual)<--(o:case)
RETURN
  f.file_name AS `File Name`,
  f.file_type AS `File Type`,
  head(labels(parent)) AS `Association`,
  f.file_description AS `Description`,
  f.file_format AS `Format`,
  f.file_size AS `Size`,
  samp.sample_id AS `Sample ID`,
  c.case_id as `Case ID`,
  i.canine_individual_id AS `Canine ID`,
  CASE WHEN s.clinical_study_designation IS NULL 
  THEN parent.clinical_study_designation 
  ELSE s.clinical_study_designation END AS `Study Code`
  
'@

$ws.Range("D2").Value2 = $cartQuery
$ws.Range("D3").Value2 = $cartQuery
$ws.Range("D4").Value2 = $cartQuery

# Match formatting used by sibling query columns (wrap-text style)
$ws.Range("D2:D4").WrapText = $true

# Row heights grew to fit the longer wrapped text
$ws.Rows.Item(2).RowHeight = 390
$ws.Rows.Item(3).RowHeight = 390
$ws.Rows.Item(4).RowHeight = 390

# Selection left behind by the editing session
$ws.Range("C14").Select()
